$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.408.80'
$ws.Range("E2").Value = '  +0.17%  '
$ws.Range("D3").Value = '2.647.86'
$ws.Range("E3").Value = '  +0.16%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '596.54'
$ws.Range("E5").Value = '  -0.24%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '158.93'
$ws.Range("E6").Value = '  +2.78%  '
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.539'
$ws.Range("E8").Value = '  -1.15%  '
$ws.Range("D9").Value = '2.646.83'
$ws.Range("E9").Value = '  +0.14%  '
$ws.Range("E10").Value = '  -2.22%  '
$ws.Range("E11").Value = '  -0.99%  '
$ws.Range("E12").Value = '  +0.32%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.352'
$ws.Range("E13").Value = '  -0.70%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '28.01'
$ws.Range("E14").Value = '  +0.19%  '
$ws.Range("D15").Value = '3.132.01'
$ws.Range("E15").Value = '  +0.27%  '
$ws.Range("E16").Value = '  -3.21%  '
$ws.Range("D17").Value = '68.320.95'
$ws.Range("E17").Value = '  -0.29%  '
$ws.Range("D18").Value = '2.639.73'
$ws.Range("E18").Value = '  -0.22%  '
$ws.Range("E19").Value = '  +1.62%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '363.39'
$ws.Range("E20").Value = '  -0.43%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.46'
$ws.Range("E21").Value = '  -0.06%  '
$ws.Range("E22").Value = '  +0.94%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.77'
$ws.Range("E23").Value = '  -2.07%  '
$ws.Range("E24").Value = '  +0.96%  '
$ws.Range("E25").Value = '  -1.01%  '
$ws.Range("E26").Value = '  +0.06%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.88'
$ws.Range("E27").Value = '  +1.09%  '
$ws.Range("E29").Value = '  -2.96%  '
$ws.Range("E30").Value = '  +0.05%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '565.69'
$ws.Range("E31").Value = '  +0.74%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.06'
$ws.Range("E32").Value = '  -0.12%  '
$ws.Range("E33").Value = '  -0.38%  '
$ws.Range("E34").Value = '  +0.38%  '
$ws.Range("E35").Value = '  +4.27%  '
$ws.Range("E36").Value = '  -1.92%  '
$ws.Range("E37").Value = '  -0.02%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '160.68'
$ws.Range("E38").Value = '  -0.40%  '
$ws.Range("E39").Value = '  +1.68%  '
$ws.Range("E40").Value = '  -1.42%  '
$ws.Range("E41").Value = '  -0.83%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.32'
$ws.Range("E42").Value = '  -0.61%  '
$ws.Range("E43").Value = '  +0.33%  '
$ws.Range("E44").Value = '  -5.52%  '
$ws.Range("E45").Value = '  +0.11%  '
$ws.Range("E46").Value = '  +1.14%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.82'
$ws.Range("E47").Value = '  +1.76%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '21.95'
$ws.Range("E48").Value = '  +0.70%  '
$ws.Range("E49").Value = '  -0.17%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0778'
$ws.Range("E50").Value = '  -0.95%  '
$ws.Range("E51").Value = '  +1.97%  '
